$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (Name, UPI ID, Ticket Code, Scanned)
$data = @(
    @("ayush",   "aj@pytm",       "ea731e88696ba970", ""),
    @("bheru",   "bheru@paypal",  "8c6d0751770bfc06", ""),
    @("dinesh",  "din@gpay",      "aa424ef629984b38", ""),
    @("shubham", "shubh@gpay",    "c6e8fcbd208e2924", ""),
    @("satyen",  "asat@razorpay", "6f80e78819a17332", ""),
    @("a1",      "a@payed",       "da9026abc665025c", ""),
    @("b2",      "b@payed",       "cfaea30e3b402902", "Yes"),
    @("c3",      "c@payed",       "4d795f2d3f2296c2", ""),
    @("d4",      "d@payed",       "7cd82693b24af5f3", "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
